# Refresh the hourly cryptocurrency snapshot (prices + 1h volume deltas)
# pulled from coinranking.com, as produced by the scheduled scraper run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new price string would otherwise be auto-coerced into a
# floating point number by the COM layer (losing trailing zeros / exact
# formatting). Force text storage for the write, then drop the explicit
# number-format style again so the cell is left exactly as before (no
# lingering "@" style on a cell that never had one).
function Set-TextValue($addr, $value) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.ClearFormats()
}

$ws.Range("D2").Value = '23.183.61'
$ws.Range("E2").Value = '  -1.22%  '
$ws.Range("D3").Value = '1.611.21'
$ws.Range("E3").Value = '  -1.25%  '
Set-TextValue "D4" '1.003'
$ws.Range("E4").Value = '  +0.44%  '
$ws.Range("E5").Value = '  +0.32%  '
Set-TextValue "D6" '302.07'
$ws.Range("E6").Value = '  -0.78%  '
Set-TextValue "D7" '0.3782'
$ws.Range("E7").Value = '  -0.11%  '
Set-TextValue "D8" '51.78'
$ws.Range("E8").Value = '  -0.49%  '
Set-TextValue "D9" '0.3530'
$ws.Range("E9").Value = '  -3.10%  '
Set-TextValue "D10" '0.08089'
$ws.Range("E10").Value = '  -0.37%  '
Set-TextValue "D11" '1.200'
$ws.Range("E11").Value = '  -2.69%  '
Set-TextValue "D12" '1.003'
$ws.Range("E12").Value = '  +0.34%  '
Set-TextValue "D13" '22.01'
$ws.Range("E13").Value = '  -3.04%  '
Set-TextValue "D14" '6.361'
$ws.Range("E14").Value = '  -3.35%  '
Set-TextValue "D15" '7.229'
$ws.Range("E15").Value = '  -0.46%  '
Set-TextValue "D16" '0.00001206'
$ws.Range("E16").Value = '  -3.39%  '
$ws.Range("D17").Value = '1.610.72'
$ws.Range("E17").Value = '  -0.84%  '
Set-TextValue "D18" '94.22'
$ws.Range("E18").Value = '  +0.43%  '
Set-TextValue "D19" '0.06913'
$ws.Range("E19").Value = '  -0.36%  '
Set-TextValue "D20" '6.500'
$ws.Range("E20").Value = '  +0.85%  '
$ws.Range("B21").Value = 'Dai'
$ws.Range("C21").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue "D21" '1.002'
$ws.Range("E21").Value = '  +0.30%  '
$ws.Range("B22").Value = 'Avalanche'
$ws.Range("C22").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextValue "D22" '17.16'
$ws.Range("E22").Value = '  -4.33%  '
Set-TextValue "D23" '12.30'
$ws.Range("E23").Value = '  -3.69%  '
$ws.Range("D24").Value = '23.182.10'
$ws.Range("E24").Value = '  -1.27%  '
Set-TextValue "D25" '2.503'
$ws.Range("E25").Value = '  +3.46%  '
Set-TextValue "D26" '3.012'
$ws.Range("E26").Value = '  -6.71%  '
Set-TextValue "D27" '20.84'
$ws.Range("E27").Value = '  -1.91%  '
Set-TextValue "D28" '150.82'
$ws.Range("E28").Value = '  +0.57%  '
Set-TextValue "D29" '5.223'
$ws.Range("E29").Value = '  -1.18%  '
Set-TextValue "D30" '132.39'
$ws.Range("E30").Value = '  -1.82%  '
$ws.Range("D31").Value = '1.797.21'
$ws.Range("E31").Value = '  -0.37%  '
Set-TextValue "D32" '1.071'
$ws.Range("E32").Value = '  +11.14%  '
Set-TextValue "D33" '6.464'
$ws.Range("E33").Value = '  -5.56%  '
$ws.Range("E34").Value = '  -8.91%  '
Set-TextValue "D35" '11.43'
$ws.Range("E35").Value = '  +3.30%  '
Set-TextValue "D36" '0.02697'
$ws.Range("E36").Value = '  -3.87%  '
Set-TextValue "D37" '0.08733'
$ws.Range("E37").Value = '  -1.30%  '
Set-TextValue "D38" '0.2445'
$ws.Range("E38").Value = '  -3.75%  '
Set-TextValue "D39" '0.06917'
$ws.Range("E39").Value = '  -4.18%  '
Set-TextValue "D40" '5.820'
$ws.Range("E40").Value = '  -4.92%  '
Set-TextValue "D41" '1.323'
$ws.Range("E41").Value = '  -2.43%  '
Set-TextValue "D42" '0.6851'
$ws.Range("E42").Value = '  -3.73%  '
Set-TextValue "D43" '11.94'
$ws.Range("E43").Value = '  -3.37%  '
Set-TextValue "D44" '15.32'
$ws.Range("E44").Value = '  -5.88%  '
$ws.Range("E45").Value = '  +0.26%  '
Set-TextValue "D46" '0.6281'
$ws.Range("E46").Value = '  -3.76%  '
Set-TextValue "D47" '3.940'
$ws.Range("E47").Value = '  -1.55%  '
Set-TextValue "D48" '2.238'
$ws.Range("E48").Value = '  -4.43%  '
Set-TextValue "D49" '0.07847'
$ws.Range("E49").Value = '  -2.09%  '
Set-TextValue "D50" '126.76'
$ws.Range("E50").Value = '  +0.78%  '
Set-TextValue "D51" '1.165'
$ws.Range("E51").Value = '  -3.56%  '
